$d = $word.ActiveDocument

# The document contains three occurrences of an <id>...</id> tag split
# across three runs: "<id>" (Courier New styled), the bare id value
# (plain black styled), and "</id>" (Courier New styled). Each trio
# needs to collapse into a single run containing the full "<id>VALUE</id>"
# text, keeping the Courier New styling of the surrounding tag runs.
# Word's Find/Replace naturally merges a multi-run match into a single
# run using the (first) run's formatting, which is exactly the desired
# end state here.

$ids = @("p098r_1", "p098r_2", "p098r_3")

foreach ($id in $ids) {
    $old = "<id>" + $id + "</id>"
    $r = $d.Content
    $found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $old, 2)
}
